$wb = $excel.ActiveWorkbook

$details = $wb.Worksheets.Item("Details")
$ws = $wb.Worksheets.Item("Sheet1")

# Rename Sheet1 -> AddProductTest
$ws.Name = "AddProductTest"

# --- Populate header row ---
$ws.Cells.Item(1,1).Value = "UserName"
$ws.Cells.Item(1,2).Value = "Password"
$ws.Cells.Item(1,3).Value = "ProductName"

# --- Populate data row ---
$ws.Cells.Item(2,1).Value = "standard_user"
$ws.Cells.Item(2,2).Value = "secret_sauce"
$ws.Cells.Item(2,3).Value = "Sauce Labs Backpack;Test.allTheThings() T-Shirt (Red);Sauce Labs Bike Light"

# --- Copy cell formatting (styles) from the Details sheet so the same
#     style indices already present in styles.xml get reused ---
$details.Range("A1:B1").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null

$details.Range("A2:B2").Copy() | Out-Null
$ws.Range("A2:B2").PasteSpecial(-4122) | Out-Null

$details.Range("A1").Copy() | Out-Null
$ws.Cells.Item(1,3).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Selections / active sheet ---
$details.Activate() | Out-Null
$details.Range("C11").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("C2").Select() | Out-Null

Write-Host "done"
